$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 24

$ws.Cells.Item($newRow, 1).Value = "Record"
$ws.Cells.Item($newRow, 2).Value = "Balanço Geral"
$ws.Cells.Item($newRow, 3).Value = "Defesa Civil"
$ws.Cells.Item($newRow, 4).Value = "2025-04-01T13:06"
$ws.Cells.Item($newRow, 5).Value = "Neutro"
$ws.Cells.Item($newRow, 6).Value = "Defesa Civil de Campos realiza demolição parcial de prédio com risco de desabar. Repórter *ao vivo*. Vídeo com depoimento do secretário da defesa Civil, Alcemir Pascoutto. "
